$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated DM_Stat (C) and P_Value (D) results from corrected Diebold-Mariano test
$ws.Range("C2").Value = 1.083054266099785
$ws.Range("D2").Value = 0.2905124201224141

$ws.Range("C3").Value = 1.160630018492204
$ws.Range("D3").Value = 0.2582322809807329

$ws.Range("C4").Value = 2.088298675658246
$ws.Range("D4").Value = 0.04855903745225643

$ws.Range("C5").Value = 3.191528809958149
$ws.Range("D5").Value = 0.004215345639827062

$ws.Range("C6").Value = 0.3534512444611433
$ws.Range("D6").Value = 0.7271146992395963

$ws.Range("C7").Value = 1.543779111424143
$ws.Range("D7").Value = 0.136905425425808

$ws.Range("C8").Value = 2.271434376891581
$ws.Range("D8").Value = 0.03325278654351571

$ws.Range("C9").Value = 0.7820604092887072
$ws.Range("D9").Value = 0.4425154555699029

$ws.Range("C10").Value = 1.957663613524952
$ws.Range("D10").Value = 0.06307346878701292
$ws.Range("G10").Value = "No"

$ws.Range("C11").Value = 1.415353912214771
$ws.Range("D11").Value = 0.170965476135875
